# Actualización automática 2025-06-01 08:00:06
#
# The monthly sales tracker rolls forward by one month:
#   - "VENTA MENSUAL" shifts its 4 month columns (C..F) one column to the
#     left (the oldest month "febrero" is dropped, a new empty month
#     "junio" appears on the right) and the header labels shift too.
#   - "VENTAS POR GRUPO" held the breakdown-by-product-category for the
#     month that is now being dropped from the rolling window (the old
#     "mayo" figures), so those per-category numbers - and the matching
#     "x de 28" counts in the totals row - reset to 0.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "VENTAS POR GRUPO": zero out the category totals that belonged
# to the month which just rolled out of the "VENTA MENSUAL" window.
# ---------------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

$grupoZeroCells = @("E2", "G2", "L2", "L3", "K4", "L8", "L10", "L13", "E17", "F17", "L17", "L21", "E24", "L25")
foreach ($ref in $grupoZeroCells) {
    $wsGrupo.Range($ref).Value2 = 0
}

# Matching "x de 28" occurrence counters on the totals row reset to 0.
$grupoCountCells = @("E30", "F30", "G30", "K30", "L30")
foreach ($ref in $grupoCountCells) {
    $wsGrupo.Range($ref).Value2 = "0 de 28"
}

# ---------------------------------------------------------------------
# Sheet "VENTA MENSUAL": shift months one column left, drop the oldest
# month, append a fresh (empty) month on the right.
# ---------------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

# Column header labels shift left; a new month label enters on the right.
$wsMensual.Range("C1").Value2 = "marzo"
$wsMensual.Range("D1").Value2 = "abril"
$wsMensual.Range("E1").Value2 = "mayo"
$wsMensual.Range("F1").Value2 = "junio"

# Data rows (2-29 detail rows, 30 totals row) shift left: new C/D/E take
# the old D/E/F values, new F is the freshly-opened (empty/zero) month.
for ($row = 2; $row -le 30; $row++) {
    $oldD = $wsMensual.Cells.Item($row, 4).Value2
    $oldE = $wsMensual.Cells.Item($row, 5).Value2
    $oldF = $wsMensual.Cells.Item($row, 6).Value2

    $wsMensual.Cells.Item($row, 3).Value2 = $oldD
    $wsMensual.Cells.Item($row, 4).Value2 = $oldE
    $wsMensual.Cells.Item($row, 5).Value2 = $oldF
    $wsMensual.Cells.Item($row, 6).Value2 = 0
}

# Column widths were auto-fit to the new content and shifted similarly.
$wsMensual.Columns.Item(3).ColumnWidth = 12.166666666666666
$wsMensual.Columns.Item(5).ColumnWidth = 13.166666666666666
$wsMensual.Columns.Item(6).ColumnWidth = 10.166666666666666
